# Update the "想去人数" (column F) counts that changed between scrapes.
# The same updates apply to both the "展览" sheet and the "全部类型" sheet,
# which mirror each other's data in this workbook.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1080
    3  = 4
    5  = 3075
    7  = 2410
    9  = 118
    11 = 1202
    14 = 7
    15 = 1081
    16 = 292
    17 = 310
    18 = 15
    19 = 19
    21 = 61
    22 = 81
    23 = 63
    24 = 8
    25 = 232
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
